{"js": "// Update Lilly Nguyen's icebreaker entry:\n//   1. Change the date line from \"9/2/2023 Lilly Nguyen\" to \"9/5/2023 Lilly Nguyen\".\n//   2. Rewrite her bio paragraph with the new wording.\n//   3. Add a new blank paragraph at the very end of the document.\n\nconst body = context.document.body;\n\n// 1) Fix the date in Lilly Nguyen's date/name line.\nconst dateResults = body.search(\"9/2/2023 Lilly Nguyen\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"9/5/2023 Lilly Nguyen\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Rewrite Lilly Nguyen's bio paragraph (locate it via a stable leading phrase).\nconst bioResults = body.search(\"Hello. I am born and raised in Jacksonville\", { matchCase: true });\nbioResults.load(\"items\");\nawait context.sync();\n\nif (bioResults.items.length > 0) {\n  const bioRange = bioResults.items[0].paragraphs.getFirst();\n  bioRange.insertText(\n    \"Hello, my name is Lilly Nguyen and I am born and raised in Jacksonville FL. \" +\n    \"I enjoy playing video games and fishing in my spare time. I enjoy playing video games such as \" +\n    \"Valorant, TFT, and BTD6. As for fishing, I mostly do saltwater fishing and have found the magical \" +\n    \"place called the Jax Pier. I hope to eventually fish an alligator. \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 3) Append a new blank paragraph after the last paragraph in the document.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Update Lilly Nguyen's icebreaker entry:\n#   1. Change the date line from \"9/2/2023 Lilly Nguyen\" to \"9/5/2023 Lilly Nguyen\".\n#   2. Rewrite her bio paragraph with the new wording.\n#   3. Add a new blank paragraph at the very end of the document.\n\n$d = $word.ActiveDocument\n\n# 1) Fix the date in Lilly Nguyen's date/name line.\n$find = $d.Content.Find\n$find.Text = \"9/2/2023 Lilly Nguyen\"\n$find.Replacement.Text = \"9/5/2023 Lilly Nguyen\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Rewrite Lilly Nguyen's bio paragraph. Done as two surgical replacements so the\n#    untouched \"Valorant\" run (and its spell-check proofErr wrapper) in the middle of the\n#    paragraph is left intact, just like in the target edit.\n$bioFind1 = $d.Content.Find\n$bioFind1.Text = \"Hello. I am born and raised in Jacksonville, FL and have loved everyday of it. I love the water and friendly alligators that are residents here. The only neighbor that I have come to despise, and hate are roaches. My hobbies are playing games and fishing. I like games such as \"\n$bioFind1.Execute() | Out-Null\n\nif ($bioFind1.Found) {\n    $bioRange1 = $d.Content.Duplicate\n    $bioRange1.Start = $bioFind1.Parent.Start\n    $bioRange1.End = $bioFind1.Parent.End\n    $bioRange1.Text = \"Hello, my name is Lilly Nguyen and I am born and raised in Jacksonville FL. I enjoy playing video games and fishing in my spare time. I enjoy playing video games such as \"\n}\n\n$bioFind2 = $d.Content.Find\n$bioFind2.Text = \", TFT, and BTD6. Fishing is a great past time for getting out of the house. I hope I can fish a gator one day as well.\"\n$bioFind2.Execute() | Out-Null\n\nif ($bioFind2.Found) {\n    $bioRange2 = $d.Content.Duplicate\n    $bioRange2.Start = $bioFind2.Parent.Start\n    $bioRange2.End = $bioFind2.Parent.End\n    $bioRange2.Text = \", TFT, and BTD6. As for fishing, I mostly do saltwater fishing and have found the magical place called the Jax Pier. I hope to eventually fish an alligator. \"\n}\n\n# 3) Append a new blank paragraph after the last paragraph in the document.\n$endRange = $d.Content\n$endRange.Collapse(0) | Out-Null\n$endRange.InsertParagraphAfter() | Out-Null\n"}
